$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 5. Индикатор (B4): update indicator wording and make it bold,
# matching the refreshed SDG 8.6.1 indicator text.
$ws.Range("B4").Value = "8.6.1  Доля молодежи (в возрасте от 15 до 24 лет), которая не учится, не работает и не приобретает профессиональных навыков "
$ws.Range("B4").Font.Bold = $true

# Сайт организации (B10): update the website and make it bold.
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B10").Font.Bold = $true
